$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 27778310
$ws.Range("I28").Value = 50000460
$ws.Range("J28").Value = 625
$ws.Range("K28").Value = 50000460
$ws.Range("L28").Value = 625
$ws.Range("M28").Value = -49999975
$ws.Range("N28").Value = -1595

$ws.Range("H43").Value = 1478.3334
$ws.Range("I43").Value = 766.6667
$ws.Range("J43").Value = 1834.1666
$ws.Range("K43").Value = 766.6667
$ws.Range("L43").Value = 1834.1666
$ws.Range("M43").Value = -697.6667
$ws.Range("N43").Value = -1972.1666

$ws.Range("H98").Value = 39529576
$ws.Range("I98").Value = 16668484
$ws.Range("J98").Value = 85251760
$ws.Range("K98").Value = 16668484
$ws.Range("L98").Value = 85251760
$ws.Range("M98").Value = -16666986
$ws.Range("N98").Value = -85254756

$ws.Range("H112").Value = 1151.5625
$ws.Range("I112").Value = 606.4286
$ws.Range("J112").Value = 1575.5555
$ws.Range("K112").Value = 1819.2858
$ws.Range("L112").Value = 4726.666499999999
$ws.Range("M112").Value = -711.2857999999999
$ws.Range("N112").Value = -6942.666499999999

$ws.Range("H122").Value = 39529576
$ws.Range("I122").Value = 16668484
$ws.Range("J122").Value = 85251760
$ws.Range("K122").Value = 50005452
$ws.Range("L122").Value = 255755280
$ws.Range("M122").Value = -50003002
$ws.Range("N122").Value = -255760180

$ws.Range("H137").Value = 39500750
$ws.Range("I137").Value = 10417990
$ws.Range("J137").Value = 89356910
$ws.Range("K137").Value = 31253970
$ws.Range("L137").Value = 268070730
$ws.Range("M137").Value = -31251420
$ws.Range("N137").Value = -268075830

$ws.Range("H138").Value = 2324.1428
$ws.Range("I138").Value = 1900.9642
$ws.Range("J138").Value = 2606.262
$ws.Range("K138").Value = 5702.892599999999
$ws.Range("L138").Value = 7818.786
$ws.Range("M138").Value = -562.8925999999992
$ws.Range("N138").Value = -18098.786

$ws.Range("H141").Value = 3747.5
$ws.Range("I141").Value = 4097
$ws.Range("K141").Value = 12291
$ws.Range("M141").Value = -7111

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 455855.9
$ws.Range("I45").Value = 527601.5600000001
$ws.Range("K45").Value = 527601.5600000001
$ws.Range("M45").Value = -527224.5600000001

$ws.Range("H61").Value = 5025466.5
$ws.Range("I61").Value = 2451988
$ws.Range("J61").Value = 19608510
$ws.Range("K61").Value = 2451988
$ws.Range("L61").Value = 19608510
$ws.Range("M61").Value = -2451776
$ws.Range("N61").Value = -19608934

$ws.Range("H74").Value = 77384100
$ws.Range("I74").Value = 84525176
$ws.Range("J74").Value = 66672476
$ws.Range("K74").Value = 84525176
$ws.Range("L74").Value = 66672476
$ws.Range("M74").Value = -84524302
$ws.Range("N74").Value = -66674224

$ws.Range("H77").Value = 77384100
$ws.Range("I77").Value = 84525176
$ws.Range("J77").Value = 66672476
$ws.Range("K77").Value = 422625880
$ws.Range("L77").Value = 333362380
$ws.Range("M77").Value = -422621512
$ws.Range("N77").Value = -333371116

$ws.Range("H88").Value = 4957.4165
$ws.Range("I88").Value = 2119.8
$ws.Range("J88").Value = 6984.2856
$ws.Range("K88").Value = 2119.8
$ws.Range("L88").Value = 6984.2856
$ws.Range("M88").Value = -1713.8
$ws.Range("N88").Value = -7796.2856

$ws.Range("H91").Value = 4957.4165
$ws.Range("I91").Value = 2119.8
$ws.Range("J91").Value = 6984.2856
$ws.Range("K91").Value = 2119.8
$ws.Range("L91").Value = 6984.2856
$ws.Range("M91").Value = -715.8000000000002
$ws.Range("N91").Value = -9792.285599999999

$ws.Range("H122").Value = 3087.3333
$ws.Range("I122").Value = 1506
$ws.Range("J122").Value = 6250
$ws.Range("K122").Value = 4518
$ws.Range("L122").Value = 18750
$ws.Range("M122").Value = -2068
$ws.Range("N122").Value = -23650

$ws.Range("H132").Value = 18524558
$ws.Range("I132").Value = 19614684
$ws.Range("J132").Value = 13891526
$ws.Range("K132").Value = 58844052
$ws.Range("L132").Value = 41674578
$ws.Range("M132").Value = -58841522
$ws.Range("N132").Value = -41679638

$ws.Range("H136").Value = 5025466.5
$ws.Range("I136").Value = 2451988
$ws.Range("J136").Value = 19608510
$ws.Range("K136").Value = 7355964
$ws.Range("L136").Value = 58825530
$ws.Range("M136").Value = -7353414
$ws.Range("N136").Value = -58830630

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1936.9395
$ws.Range("I86").Value = 1951.0416
$ws.Range("J86").Value = 1485.6666
$ws.Range("K86").Value = 1951.0416
$ws.Range("L86").Value = 1485.6666
$ws.Range("M86").Value = -828.0416
$ws.Range("N86").Value = -3731.6666

$ws.Range("H89").Value = 1936.9395
$ws.Range("I89").Value = 1951.0416
$ws.Range("J89").Value = 1485.6666
$ws.Range("K89").Value = 9755.208000000001
$ws.Range("L89").Value = 7428.333000000001
$ws.Range("M89").Value = -4139.208000000001
$ws.Range("N89").Value = -18660.333

$ws.Range("H105").Value = 2184.8235
$ws.Range("I105").Value = 1968.375
$ws.Range("J105").Value = 2377.2222
$ws.Range("K105").Value = 1968.375
$ws.Range("L105").Value = 2377.2222
$ws.Range("M105").Value = -221.375
$ws.Range("N105").Value = -5871.2222

$ws.Range("H115").Value = 30684
$ws.Range("J115").Value = 30684
$ws.Range("L115").Value = 30684
$ws.Range("N115").Value = -33818

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 55568204
$ws.Range("I99").Value = 90916970
$ws.Range("J99").Value = 20142.285
$ws.Range("K99").Value = 90916970
$ws.Range("L99").Value = 20142.285
$ws.Range("M99").Value = -90915472
$ws.Range("N99").Value = -23138.285

$ws.Range("H126").Value = 55568204
$ws.Range("I126").Value = 90916970
$ws.Range("J126").Value = 20142.285
$ws.Range("K126").Value = 272750910
$ws.Range("L126").Value = 60426.855
$ws.Range("M126").Value = -272748440
$ws.Range("N126").Value = -65366.855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2565350.2
$ws.Range("J5").Value = 3032298.2
$ws.Range("L5").Value = 9096894.600000001
$ws.Range("N5").Value = -9097118.600000001

$ws.Range("H11").Value = 134.44444
$ws.Range("I11").Value = 70
$ws.Range("K11").Value = 210
$ws.Range("M11").Value = -70

$ws.Range("H22").Value = 1407.6923
$ws.Range("J22").Value = 1483.3334
$ws.Range("L22").Value = 4450.0002
$ws.Range("N22").Value = -4788.0002

$ws.Range("H27").Value = 1407.6923
$ws.Range("J27").Value = 1483.3334
$ws.Range("L27").Value = 4450.0002
$ws.Range("N27").Value = -4654.0002

$ws.Range("H68").Value = 13456
$ws.Range("I68").Value = 925
$ws.Range("J68").Value = 17633
$ws.Range("K68").Value = 2775
$ws.Range("L68").Value = 52899
$ws.Range("M68").Value = -1964
$ws.Range("N68").Value = -54521

$ws.Range("H71").Value = 13456
$ws.Range("I71").Value = 925
$ws.Range("J71").Value = 17633
$ws.Range("K71").Value = 8325
$ws.Range("L71").Value = 158697
$ws.Range("M71").Value = -4269
$ws.Range("N71").Value = -166809

$ws.Range("H117").Value = 2458.1428
$ws.Range("J117").Value = 2458.1428
$ws.Range("L117").Value = 7374.428400000001
$ws.Range("N117").Value = -14258.4284

$ws.Range("H129").Value = 2043.7333
$ws.Range("I129").Value = 1516.25
$ws.Range("K129").Value = 4548.75
$ws.Range("M129").Value = 451.25

$ws.Range("H135").Value = 2565350.2
$ws.Range("J135").Value = 3032298.2
$ws.Range("L135").Value = 27290683.8
$ws.Range("N135").Value = -27295753.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2406910
$ws.Range("I70").Value = 1548778.4
$ws.Range("J70").Value = 3642619.5
$ws.Range("K70").Value = 1548778.4
$ws.Range("L70").Value = 3642619.5
$ws.Range("M70").Value = -1548508.4
$ws.Range("N70").Value = -3643159.5

$ws.Range("H73").Value = 2406910
$ws.Range("I73").Value = 1548778.4
$ws.Range("J73").Value = 3642619.5
$ws.Range("K73").Value = 1548778.4
$ws.Range("L73").Value = 3642619.5
$ws.Range("M73").Value = -1547842.4
$ws.Range("N73").Value = -3644491.5

$ws.Range("H107").Value = 349.85715
$ws.Range("I107").Value = 90.90000000000001
$ws.Range("J107").Value = 997.25
$ws.Range("K107").Value = 90.90000000000001
$ws.Range("L107").Value = 997.25
$ws.Range("M107").Value = 1829.1
$ws.Range("N107").Value = -4837.25

$ws.Range("H132").Value = 14107841
$ws.Range("I132").Value = 16509407
$ws.Range("J132").Value = 10105231
$ws.Range("K132").Value = 49528221
$ws.Range("L132").Value = 30315693
$ws.Range("M132").Value = -49525691
$ws.Range("N132").Value = -30320753

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7274.5713
$ws.Range("I22").Value = 740
$ws.Range("J22").Value = 23611
$ws.Range("K22").Value = 740
$ws.Range("L22").Value = 23611
$ws.Range("M22").Value = -445
$ws.Range("N22").Value = -24201

$ws.Range("H27").Value = 7274.5713
$ws.Range("I27").Value = 740
$ws.Range("J27").Value = 23611
$ws.Range("K27").Value = 740
$ws.Range("L27").Value = 23611
$ws.Range("M27").Value = -633
$ws.Range("N27").Value = -23825

$ws.Range("H64").Value = 221028
$ws.Range("J64").Value = 267535
$ws.Range("L64").Value = 267535
$ws.Range("N64").Value = -267985

$ws.Range("H67").Value = 221028
$ws.Range("J67").Value = 267535
$ws.Range("L67").Value = 267535
$ws.Range("N67").Value = -269095

$ws.Range("H132").Value = 1669375.4
$ws.Range("I132").Value = 2300107
$ws.Range("J132").Value = 6537.1816
$ws.Range("K132").Value = 6900321
$ws.Range("L132").Value = 19611.5448
$ws.Range("M132").Value = -6897791
$ws.Range("N132").Value = -24671.5448

$ws.Range("H136").Value = 2029692.9
$ws.Range("I136").Value = 3269171.5
$ws.Range("J136").Value = 1455.4546
$ws.Range("K136").Value = 9807514.5
$ws.Range("L136").Value = 4366.3638
$ws.Range("M136").Value = -9804964.5
$ws.Range("N136").Value = -9466.363799999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 984.15625
$ws.Range("I122").Value = 833.46155
$ws.Range("J122").Value = 1637.1666
$ws.Range("K122").Value = 2500.38465
$ws.Range("L122").Value = 4911.4998
$ws.Range("M122").Value = -50.38464999999997
$ws.Range("N122").Value = -9811.4998

$ws.Range("H136").Value = 11918.81
$ws.Range("I136").Value = 8572
$ws.Range("J136").Value = 20285.834
$ws.Range("K136").Value = 25716
$ws.Range("L136").Value = 60857.50199999999
$ws.Range("M136").Value = -65957.50199999999
